$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D2 value to 2
$ws.Range("D2").Value = 2

# Set D3 formula to sum D2:D2
$ws.Range("D3").Formula = "=SUM(D2:D2)"

# Update the active selection to D4
$ws.Range("D4").Select()
